# Re-apply the header-cell character formatting ("Month" / "Savings") that
# the Apache POI 4.1.0 -> 5.2.3 upgrade re-serializes from w:val="true"/
# "false" to w:val="on"/"off". Semantically the runs stay Bold / not-Italic /
# not-Strike — we just re-assert that formatting through the Word object
# model so the rPr is rewritten in the new writer's canonical form.

$d = $word.ActiveDocument

$table = $d.Tables.Item(1)

$headerCells = @($table.Cell(1, 1).Range, $table.Cell(1, 2).Range)

foreach ($cellRange in $headerCells) {
    $cellRange.Font.Bold = $true
    $cellRange.Font.Italic = $false
    $cellRange.Font.StrikeThrough = $false
}
